$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.708.48"
$ws.Range("E2").Value = "  +2.24%  "
$ws.Range("D3").Value = "1.871.40"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").Value = "'323.80"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "'0.4584"
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("D8").Value = "'0.3852"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "'0.07854"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "'0.9920"
$ws.Range("E10").Value = "  +3.24%  "
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("D12").Value = "1.884.32"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "'6.962"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").Value = "'5.688"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "'0.06964"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "'88.30"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "'1.006"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "'0.00001003"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").Value = "'16.80"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "28.746.60"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").Value = "'5.274"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "'11.00"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'2.128"
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("D25").Value = "2.130.92"
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("D26").Value = "'153.33"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("D27").Value = "'19.21"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D28").Value = "'5.761"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").Value = "'1.946"
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("D30").Value = "'118.77"
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("D31").Value = "'0.09304"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "'0.9146"
$ws.Range("E32").Value = "  -2.82%  "
$ws.Range("D33").Value = "'5.289"
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("D34").Value = "'1.333"
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("D35").Value = "'3.307"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").Value = "'0.05742"
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("D38").Value = "'0.02075"
$ws.Range("E38").Value = "  -1.80%  "
$ws.Range("D39").Value = "'7.698"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").Value = "'0.5632"
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("D41").Value = "'0.1786"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").Value = "'9.843"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("D43").Value = "'0.07187"
$ws.Range("E43").Value = "  -2.29%  "
$ws.Range("D44").Value = "'11.79"
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D45").Value = "'0.5267"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'2.122"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").Value = "'1.119"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("D48").Value = "'1.824"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").Value = "'113.27"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("E50").Value = "  +4.07%  "
$ws.Range("D51").Value = "'1.004"
